# Apply weekly update to Fruta / hortaliza (Kiwi) data sheet.
# The data rows (2..14) had their values permuted: the row that now
# appears at a given row index previously appeared at a different row
# index. Columns D, L, M, N, O, P, Q, S, T carry the data that moves;
# columns A, B, C, E, F, G, H, I, J, K, R stay constant across rows so
# they do not need to change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping: destination row -> source row (values to copy from source row
# in the ORIGINAL workbook into the destination row).
$rowMap = @{
    2  = 3
    3  = 7
    4  = 5
    5  = 11
    6  = 2
    7  = 8
    8  = 13
    9  = 10
    10 = 14
    11 = 12
    12 = 6
    13 = 4
    14 = 9
}

# Columns that carry the per-row data which gets permuted.
$cols = @("D", "L", "M", "N", "O", "P", "Q", "S", "T")

# Snapshot all original values for the columns/rows involved before
# writing anything, since several rows both provide and receive values.
$original = @{}
foreach ($col in $cols) {
    for ($r = 2; $r -le 14; $r++) {
        $original["$col$r"] = $ws.Range("$col$r").Value2
    }
}

# Write the new values based on the row mapping.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $original["$col$srcRow"]
    }
}
